$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Cells.Item(2, 4).Value = '64.637.64'
$ws.Cells.Item(2, 5).Value = '  -1.34%  '
$ws.Cells.Item(3, 4).Value = '3.514.70'
$ws.Cells.Item(3, 5).Value = '  -2.04%  '
$ws.Cells.Item(4, 5).Value = '  +0.08%  '
$ws.Cells.Item(5, 4).NumberFormat = '@'
$ws.Cells.Item(5, 4).Value = '586.69'
$ws.Cells.Item(5, 5).Value = '  -2.93%  '
$ws.Cells.Item(6, 4).NumberFormat = '@'
$ws.Cells.Item(6, 4).Value = '132.68'
$ws.Cells.Item(6, 5).Value = '  -2.28%  '
$ws.Cells.Item(7, 4).Value = '3.515.21'
$ws.Cells.Item(7, 5).Value = '  -2.04%  '
$ws.Cells.Item(8, 5).Value = '  +0.03%  '
$ws.Cells.Item(9, 4).NumberFormat = '@'
$ws.Cells.Item(9, 4).Value = '0.491'
$ws.Cells.Item(9, 5).Value = '  -0.94%  '
$ws.Cells.Item(10, 5).Value = '  -0.34%  '
$ws.Cells.Item(11, 5).Value = '  +0.91%  '
$ws.Cells.Item(12, 4).NumberFormat = '@'
$ws.Cells.Item(12, 4).Value = '0.388'
$ws.Cells.Item(12, 5).Value = '  -0.97%  '
$ws.Cells.Item(13, 4).Value = '4.109.79'
$ws.Cells.Item(13, 5).Value = '  -1.85%  '
$ws.Cells.Item(14, 4).NumberFormat = '@'
$ws.Cells.Item(14, 4).Value = '27.75'
$ws.Cells.Item(14, 5).Value = '  +0.06%  '
$ws.Cells.Item(15, 5).Value = '  -3.14%  '
$ws.Cells.Item(16, 5).Value = '  +0.59%  '
$ws.Cells.Item(17, 4).Value = '3.511.42'
$ws.Cells.Item(17, 5).Value = '  -1.87%  '
$ws.Cells.Item(18, 4).Value = '64.617.29'
$ws.Cells.Item(18, 5).Value = '  -1.44%  '
$ws.Cells.Item(19, 4).NumberFormat = '@'
$ws.Cells.Item(19, 4).Value = '9.98'
$ws.Cells.Item(19, 5).Value = '  -0.66%  '
$ws.Cells.Item(20, 4).NumberFormat = '@'
$ws.Cells.Item(20, 4).Value = '14.20'
$ws.Cells.Item(20, 5).Value = '  -2.76%  '
$ws.Cells.Item(21, 5).Value = '  -3.95%  '
$ws.Cells.Item(22, 4).NumberFormat = '@'
$ws.Cells.Item(22, 4).Value = '391.50'
$ws.Cells.Item(22, 5).Value = '  -0.91%  '
$ws.Cells.Item(23, 4).NumberFormat = '@'
$ws.Cells.Item(23, 4).Value = '0.579'
$ws.Cells.Item(23, 5).Value = '  -1.44%  '
$ws.Cells.Item(24, 4).Value = '3.654.68'
$ws.Cells.Item(24, 5).Value = '  -2.00%  '
$ws.Cells.Item(25, 4).NumberFormat = '@'
$ws.Cells.Item(25, 4).Value = '73.86'
$ws.Cells.Item(25, 5).Value = '  -0.65%  '
$ws.Cells.Item(26, 5).Value = '  +0.14%  '
$ws.Cells.Item(27, 4).NumberFormat = '@'
$ws.Cells.Item(27, 4).Value = '0.0000111'
$ws.Cells.Item(27, 5).Value = '  -4.76%  '
$ws.Cells.Item(28, 4).NumberFormat = '@'
$ws.Cells.Item(28, 4).Value = '1.58'
$ws.Cells.Item(28, 5).Value = '  -5.16%  '
$ws.Cells.Item(29, 4).NumberFormat = '@'
$ws.Cells.Item(29, 4).Value = '7.46'
$ws.Cells.Item(29, 5).Value = '  -8.18%  '
$ws.Cells.Item(30, 4).NumberFormat = '@'
$ws.Cells.Item(30, 4).Value = '0.999'
$ws.Cells.Item(30, 5).Value = '  -0.15%  '
$ws.Cells.Item(31, 4).NumberFormat = '@'
$ws.Cells.Item(31, 4).Value = '2.27'
$ws.Cells.Item(31, 5).Value = '  -5.29%  '
$ws.Cells.Item(32, 5).Value = '  -5.44%  '
$ws.Cells.Item(33, 4).Value = '3.519.04'
$ws.Cells.Item(33, 5).Value = '  -1.72%  '
$ws.Cells.Item(35, 4).NumberFormat = '@'
$ws.Cells.Item(35, 4).Value = '24.02'
$ws.Cells.Item(35, 5).Value = '  -1.85%  '
$ws.Cells.Item(36, 5).Value = '  -0.73%  '
$ws.Cells.Item(37, 4).NumberFormat = '@'
$ws.Cells.Item(37, 4).Value = '5.28'
$ws.Cells.Item(37, 5).Value = '  -0.42%  '
$ws.Cells.Item(38, 4).NumberFormat = '@'
$ws.Cells.Item(38, 4).Value = '1.60'
$ws.Cells.Item(38, 5).Value = '  +0.18%  '
$ws.Cells.Item(39, 4).NumberFormat = '@'
$ws.Cells.Item(39, 4).Value = '171.03'
$ws.Cells.Item(39, 5).Value = '  +0.00%  '
$ws.Cells.Item(40, 4).NumberFormat = '@'
$ws.Cells.Item(40, 4).Value = '7.00'
$ws.Cells.Item(40, 5).Value = '  -1.23%  '
$ws.Cells.Item(41, 5).Value = '  -2.91%  '
$ws.Cells.Item(42, 2).Value = 'EnergySwap'
$ws.Cells.Item(42, 3).Value = 'https://coinranking.com/coin/SbWqqTui-+energyswap-ens'
$ws.Cells.Item(42, 4).NumberFormat = '@'
$ws.Cells.Item(42, 4).Value = '26.65'
$ws.Cells.Item(42, 5).Value = '  +1.88%  '
$ws.Cells.Item(43, 2).Value = 'Mantle'
$ws.Cells.Item(43, 3).Value = 'https://coinranking.com/coin/BoI4ux0nd+mantle-mnt'
$ws.Cells.Item(43, 4).NumberFormat = '@'
$ws.Cells.Item(43, 4).Value = '0.814'
$ws.Cells.Item(43, 5).Value = '  -3.78%  '
$ws.Cells.Item(44, 5).Value = '  +0.08%  '
$ws.Cells.Item(45, 4).NumberFormat = '@'
$ws.Cells.Item(45, 4).Value = '42.15'
$ws.Cells.Item(45, 5).Value = '  -2.75%  '
$ws.Cells.Item(46, 4).NumberFormat = '@'
$ws.Cells.Item(46, 4).Value = '1.22'
$ws.Cells.Item(46, 5).Value = '  -1.96%  '
$ws.Cells.Item(47, 5).Value = '  -2.76%  '
$ws.Cells.Item(48, 4).NumberFormat = '@'
$ws.Cells.Item(48, 4).Value = '1.66'
$ws.Cells.Item(48, 5).Value = '  -2.19%  '
$ws.Cells.Item(49, 4).Value = '2.453.37'
$ws.Cells.Item(49, 5).Value = '  -0.23%  '
$ws.Cells.Item(50, 4).NumberFormat = '@'
$ws.Cells.Item(50, 4).Value = '6.89'
$ws.Cells.Item(50, 5).Value = '  -2.22%  '
$ws.Cells.Item(51, 4).NumberFormat = '@'
$ws.Cells.Item(51, 4).Value = '0.904'
$ws.Cells.Item(51, 5).Value = '  +2.82%  '
